# Scheduled runner update: refresh computed market-profit columns (H:N)
# on the Leve profit sheets. Values come from an external price-data sync;
# columns are plain numbers (no formulas) so we just overwrite the cells.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 107
$ws.Range("H107").Value = 900
$ws.Range("I107").Value = 993.75
$ws.Range("K107").Value = 993.75
$ws.Range("M107").Value = 926.25

# Row 111
$ws.Range("H111").Value = 3916.2
$ws.Range("I111").Value = 4391.375
$ws.Range("J111").Value = 2015.5
$ws.Range("K111").Value = 13174.125
$ws.Range("L111").Value = 6046.5
$ws.Range("M111").Value = -10107.125
$ws.Range("N111").Value = -12180.5

# Row 132
$ws.Range("H132").Value = 2275.7896
$ws.Range("I132").Value = 2310.2703
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 6930.8109
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -4400.8109
$ws.Range("N132").Value = -8060

# Row 135
$ws.Range("H135").Value = 20006282
$ws.Range("I135").Value = 601.65
$ws.Range("J135").Value = 100029010
$ws.Range("K135").Value = 5414.849999999999
$ws.Range("L135").Value = 900261090
$ws.Range("M135").Value = -2879.849999999999
$ws.Range("N135").Value = -900266160

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 26890.064
$ws.Range("I32").Value = 29088.838
$ws.Range("J32").Value = 3253.25
$ws.Range("K32").Value = 29088.838
$ws.Range("L32").Value = 3253.25
$ws.Range("M32").Value = -28801.838
$ws.Range("N32").Value = -3827.25

# Row 45
$ws.Range("H45").Value = 2940.4285
$ws.Range("I45").Value = 2664.4546
$ws.Range("J45").Value = 3952.3333
$ws.Range("K45").Value = 2664.4546
$ws.Range("L45").Value = 3952.3333
$ws.Range("M45").Value = -2287.4546
$ws.Range("N45").Value = -4706.3333

# Row 74
$ws.Range("H74").Value = 111112160
$ws.Range("I74").Value = 111112160
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 111112160
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -111111286
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 111112160
$ws.Range("I77").Value = 111112160
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 555560800
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -555556432
$ws.Range("N77").ClearContents()

# Row 122
$ws.Range("H122").Value = 3668.4443
$ws.Range("I122").Value = 2800.8
$ws.Range("K122").Value = 8402.400000000001
$ws.Range("M122").Value = -5952.400000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 1044.4688
$ws.Range("I80").Value = 702.6667
$ws.Range("J80").Value = 1249.55
$ws.Range("K80").Value = 702.6667
$ws.Range("L80").Value = 1249.55
$ws.Range("M80").Value = 295.3333
$ws.Range("N80").Value = -3245.55

# Row 83
$ws.Range("H83").Value = 1044.4688
$ws.Range("I83").Value = 702.6667
$ws.Range("J83").Value = 1249.55
$ws.Range("K83").Value = 3513.3335
$ws.Range("L83").Value = 6247.75
$ws.Range("M83").Value = 1478.6665
$ws.Range("N83").Value = -16231.75

# Row 107
$ws.Range("H107").Value = 887.2963
$ws.Range("I107").Value = 702
$ws.Range("K107").Value = 702
$ws.Range("M107").Value = 1218

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 15864.08
$ws.Range("I31").Value = 21128.25
$ws.Range("J31").Value = 6505.5557
$ws.Range("K31").Value = 21128.25
$ws.Range("L31").Value = 6505.5557
$ws.Range("M31").Value = -20833.25
$ws.Range("N31").Value = -7095.5557

# Row 34
$ws.Range("H34").Value = 15864.08
$ws.Range("I34").Value = 21128.25
$ws.Range("J34").Value = 6505.5557
$ws.Range("K34").Value = 21128.25
$ws.Range("L34").Value = 6505.5557
$ws.Range("M34").Value = -20926.25
$ws.Range("N34").Value = -6909.5557

# Row 107
$ws.Range("H107").Value = 1140.9375
$ws.Range("I107").Value = 400
$ws.Range("J107").Value = 1311.9231
$ws.Range("K107").Value = 400
$ws.Range("L107").Value = 1311.9231
$ws.Range("M107").Value = 1520
$ws.Range("N107").Value = -5151.9231

# Row 132
$ws.Range("H132").Value = 10571.607
$ws.Range("I132").Value = 11278.694
$ws.Range("J132").Value = 5622
$ws.Range("K132").Value = 33836.08199999999
$ws.Range("L132").Value = 16866
$ws.Range("M132").Value = -31306.08199999999
$ws.Range("N132").Value = -21926

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 2932
$ws.Range("I3").Value = 1217.4445
$ws.Range("J3").Value = 4646.5557
$ws.Range("K3").Value = 3652.3335
$ws.Range("L3").Value = 13939.6671
$ws.Range("M3").Value = -3540.3335
$ws.Range("N3").Value = -14163.6671

# Row 130
$ws.Range("H130").Value = 2249.5
$ws.Range("I130").Value = 1500
$ws.Range("J130").Value = 2999
$ws.Range("K130").Value = 4500
$ws.Range("L130").Value = 8997
$ws.Range("M130").Value = 520
$ws.Range("N130").Value = -19037

# Row 131
$ws.Range("H131").Value = 753.16
$ws.Range("J131").Value = 805.4607
$ws.Range("L131").Value = 2416.3821
$ws.Range("N131").Value = -12496.3821

# Row 138
$ws.Range("H138").Value = 144139.05
$ws.Range("I138").Value = 1266
$ws.Range("K138").Value = 3798
$ws.Range("M138").Value = 1342

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 5917407
$ws.Range("I107").Value = 326.625
$ws.Range("J107").Value = 15384736
$ws.Range("K107").Value = 326.625
$ws.Range("L107").Value = 15384736
$ws.Range("M107").Value = 1593.375
$ws.Range("N107").Value = -15388576

# Row 122
$ws.Range("H122").Value = 148148910
$ws.Range("I122").Value = 41667532
$ws.Range("K122").Value = 125002596
$ws.Range("M122").Value = -125000146

# Row 126
$ws.Range("H126").Value = 6891.304
$ws.Range("I126").Value = 6468.75
$ws.Range("J126").Value = 7857.143
$ws.Range("K126").Value = 19406.25
$ws.Range("L126").Value = 23571.429
$ws.Range("M126").Value = -16936.25
$ws.Range("N126").Value = -28511.429

# Row 132
$ws.Range("H132").Value = 51777.97
$ws.Range("I132").Value = 49228.137
$ws.Range("J132").Value = 58010.89
$ws.Range("K132").Value = 147684.411
$ws.Range("L132").Value = 174032.67
$ws.Range("M132").Value = -145154.411
$ws.Range("N132").Value = -179092.67

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 7584.4443
$ws.Range("I61").Value = 3876.6667
$ws.Range("K61").Value = 3876.6667
$ws.Range("M61").Value = -3674.6667

# Row 93
$ws.Range("H93").Value = 887.875
$ws.Range("I93").Value = 775.5
$ws.Range("K93").Value = 775.5
$ws.Range("M93").Value = 472.5

# Row 113
$ws.Range("H113").Value = 7584.4443
$ws.Range("I113").Value = 3876.6667
$ws.Range("K113").Value = 3876.6667
$ws.Range("M113").Value = -1706.6667

# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# Row 136
$ws.Range("H136").Value = 24019.5
$ws.Range("I136").Value = 26046.2
$ws.Range("J136").Value = 3752.5
$ws.Range("K136").Value = 78138.60000000001
$ws.Range("L136").Value = 11257.5
$ws.Range("M136").Value = -75588.60000000001
$ws.Range("N136").Value = -16357.5

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 3863604.5
$ws.Range("I113").Value = 4375
$ws.Range("J113").Value = 9009244
$ws.Range("K113").Value = 13125
$ws.Range("L113").Value = 27027732
$ws.Range("M113").Value = -10955
$ws.Range("N113").Value = -27032072
